# Daily data refresh of the "Pais" (countries) COVID stats table.
# - A handful of countries swap rank position versus their neighbour
#   (their case counts grew/shrank enough to overtake one another), which
#   shows up as column A text swapping between two adjacent rows.
# - A batch of rows receive updated Casos totales / Nuevos casos /
#   Casos activos / Recuperados / Muertes hoy / Muertes figures.
# - The "Datos actualizados..." footer timestamp moves from 14:05 to 15:05.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (column A) caused by re-ranking in the data refresh ---
$ws.Range("A153").Value = "Liberia"
$ws.Range("A154").Value = "Suazilandia"

$ws.Range("A157").Value = "Benin"
$ws.Range("A158").Value = "Birmania"

$ws.Range("A193").Value = "Namibia"
$ws.Range("A194").Value = "Timor Oriental"

$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("A201").Value = "Belice"

$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("A214").Value = "Islas Virgenes Britanicas"

# --- Updated case statistics (columns B:H) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1818983
$ws.Range("C4").Value = 2163
$ws.Range("D4").Value = 535361
$ws.Range("E4").Value = 1178011
$ws.Range("G4").Value = 54
$ws.Range("H4").Value = 105611

# Row 5 - Brasil
$ws.Range("B5").Value = 501985
$ws.Range("C5").Value = 3545
$ws.Range("E5").Value = 267742
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = 28872

# Row 12 - India
$ws.Range("B12").Value = 183008
$ws.Range("C12").Value = 1181
$ws.Range("D12").Value = 87179
$ws.Range("E12").Value = 90641

# Row 19 - Chile
$ws.Range("B19").Value = 85261
$ws.Range("C19").Value = 1877
$ws.Range("D19").Value = 62442
$ws.Range("E19").Value = 22316
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 503

# Row 26 - Bielorrusia
$ws.Range("B26").Value = 42556
$ws.Range("C26").Value = 898
$ws.Range("D26").Value = 18514
$ws.Range("E26").Value = 23807
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 235

# Row 42 - Camerun
$ws.Range("B42").Value = 18086
$ws.Range("C42").Value = 862
$ws.Range("D42").Value = 3909
$ws.Range("E42").Value = 13220
$ws.Range("G42").Value = 7
$ws.Range("H42").Value = 957

# Row 53
$ws.Range("B53").Value = 11412
$ws.Range("C53").Value = 31
$ws.Range("D53").Value = 6698
$ws.Range("E53").Value = 4471
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 243

# Row 79
$ws.Range("B79").Value = 3583
$ws.Range("C79").Value = 37
$ws.Range("D79").Value = 2837
$ws.Range("E79").Value = 731
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 15

# Row 110
$ws.Range("B110").Value = 1220
$ws.Range("C110").Value = 29
$ws.Range("D110").Value = 711
$ws.Range("E110").Value = 482
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 27

# Row 153 - Liberia (after swap)
$ws.Range("B153").Value = 288
$ws.Range("C153").Value = 8
$ws.Range("D153").Value = 157
$ws.Range("E153").Value = 104
$ws.Range("H153").Value = 27

# Row 154 - Suazilandia (after swap)
$ws.Range("B154").Value = 283
$ws.Range("D154").Value = 168
$ws.Range("E154").Value = 113
$ws.Range("H154").Value = 2

# Row 157 - Benin (after swap)
$ws.Range("B157").Value = 232
$ws.Range("C157").Value = 8
$ws.Range("D157").Value = 143
$ws.Range("E157").Value = 86
$ws.Range("H157").Value = 3

# Row 158 - Birmania (after swap)
$ws.Range("D158").Value = 130
$ws.Range("E158").Value = 88
$ws.Range("H158").Value = 6

# Row 193 - Namibia (after swap)
$ws.Range("C193").Value = 1
$ws.Range("D193").Value = 14
$ws.Range("E193").Value = 10

# Row 194 - Timor Oriental (after swap)
$ws.Range("B194").Value = 24
$ws.Range("D194").Value = 24
$ws.Range("E194").Value = 0

# Row 200 - Santa Lucia (after swap)
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

# Row 201 - Belice (after swap)
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# Row 213 - Papua Nueva Guinea (after swap)
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214 - Islas Virgenes Britanicas (after swap)
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1

# --- Timestamp footer update ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 15:05"
